$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: clear B2/C2 content, shift "Document the TileMap..." text into B2,
# leave A2 blank (style s="2" is preserved automatically since we only clear value)
$ws.Range("A2").Value = $null
$ws.Range("B2").Value = "Document the TileMap, b2world, axis orientations"
$ws.Range("C2").Value = $null

# Row 3: *Determine Texture Sizes moves to A3, new text "Research box 2d shapes" in B3
$ws.Range("A3").Value = "*Determine Texture Sizes"
$ws.Range("B3").Value = "Research box 2d shapes"

# Row 4: "Figure out level file format" stays in A4, new text in B4
$ws.Range("A4").Value = "Figure out level file format"
$ws.Range("B4").Value = "Figure out format for entities that have only one texture"

# Update the active selection to B4 as per the sheetView change
$ws.Range("B4").Select()
